# Applies the "WOT/mbox dataset" update: inserts a new column C (local
# predicate name) for a handful of rows in the rdf-dereferencer sheet and
# pushes the previous column C (human readable label) into column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rdf-dereferencer")

# rowNum -> [newColCValue, oldColCValue-that-moves-to-D]
$updates = [ordered]@{
    33 = @("hasKey",         "has Key")
    34 = @("pubkeyAddress",  "Address")
    35 = @("fingerprint",    "Fingerprint")
    36 = @("identity",       "Identity")
    37 = @("assurance",      "Assurance")
    38 = @("signed",         "Signed")
    39 = @("signer",         "Signer")
    60 = @("mbox",           "personal mailbox")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $newC = $pair[0]
    $oldC = $pair[1]

    # Move the previous column C text into the new column D first …
    $ws.Cells.Item($row, 4).Value = $oldC
    # … then overwrite column C with the new local-name value.
    $ws.Cells.Item($row, 3).Value = $newC
}

# Update the sheet's view/selection the same way the original author left it
# (scrolled so row 28 is at the top, with E39 selected).
$ws.Activate()
$window = $excel.ActiveWindow
$window.ScrollRow = 28
$window.ScrollColumn = 1
$ws.Range("E39").Select()
